$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 26.38884429336736
$ws.Range("H2").Value = 27.25445264277199
$ws.Range("I2").Value = 7.178023577499991
$ws.Range("J2").Value = 7.725743978589128
$ws.Range("N2").Value = 0.2480225622313161
$ws.Range("F3").Value = 27.67947691269417
$ws.Range("H3").Value = 27.74069953174611
$ws.Range("I3").Value = 8.659041378755036
$ws.Range("J3").Value = 7.913223349077937
$ws.Range("L3").Value = 3.100467703943578
$ws.Range("N3").Value = 0.8682714849357696
$ws.Range("F4").Value = 28.27288603096032
$ws.Range("H4").Value = 29.59199753506476
$ws.Range("I4").Value = 9.508121753665346
$ws.Range("J4").Value = 7.958607729535689
$ws.Range("L4").Value = 3.911977283124905
$ws.Range("N4").Value = 0.2737584616091717
$ws.Range("F5").Value = 26.84805945593586
$ws.Range("H5").Value = 27.52466958217565
$ws.Range("I5").Value = 8.140651392009939
$ws.Range("J5").Value = 7.548134475631093
$ws.Range("L5").Value = 2.712482622149522
$ws.Range("N5").Value = 1.089193391690594
$ws.Range("F6").Value = 25.54199208562319
$ws.Range("H6").Value = 26.01155197083243
$ws.Range("I6").Value = 6.843499207440418
$ws.Range("J6").Value = 7.278784620126929
$ws.Range("L6").Value = 1.587609864603051
$ws.Range("N6").Value = 0.2295687001863469
$ws.Range("F7").Value = 24.99562737902271
$ws.Range("H7").Value = 25.26681073416658
$ws.Range("I7").Value = 6.649759785242955
$ws.Range("J7").Value = 7.175367625362059
$ws.Range("L7").Value = 1.231454032230362
$ws.Range("N7").Value = 0.60217061690214
$ws.Range("F8").Value = 23.85366320880505
$ws.Range("H8").Value = 25.7606110372136
$ws.Range("I8").Value = 5.741283396547754
$ws.Range("J8").Value = 7.020936897864535
$ws.Range("L8").Value = 0.4102107802364227
$ws.Range("N8").Value = 0.3722568210927845
$ws.Range("F9").Value = 22.41577226909466
$ws.Range("H9").Value = 23.55004044172754
$ws.Range("I9").Value = 4.207604559957425
$ws.Range("J9").Value = 6.851594345650664
$ws.Range("L9").Value = 0.4220388785950696
$ws.Range("N9").Value = 0.3592497067836414
$ws.Range("F10").Value = 22.66351699941429
$ws.Range("H10").Value = 23.66748617903899
$ws.Range("I10").Value = 5.086482993125527
$ws.Range("J10").Value = 6.909275675652921
$ws.Range("L10").Value = 0.2645509247625299
$ws.Range("N10").Value = 0.9581768297437723
$ws.Range("F11").Value = 20.81667481561695
$ws.Range("H11").Value = 22.21365562209478
$ws.Range("I11").Value = 3.387440512298052
$ws.Range("J11").Value = 6.716010493195187
$ws.Range("L11").Value = 1.736199019038345
$ws.Range("N11").Value = 1.63227075919059
$ws.Range("F12").Value = 18.20863591133639
$ws.Range("H12").Value = 19.28450161634704
$ws.Range("I12").Value = 1.582492011778716
$ws.Range("J12").Value = 6.447397696910457
$ws.Range("L12").Value = 4.795408174083222
$ws.Range("N12").Value = 0.9166889945213713
$ws.Range("F13").Value = 15.62418515037838
$ws.Range("H13").Value = 16.52818048991558
$ws.Range("I13").Value = 1.091220449544552
$ws.Range("J13").Value = 6.194999230802507
$ws.Range("L13").Value = 8.014703939978563
$ws.Range("N13").Value = 0.6314215351574861
$ws.Range("F14").Value = 2.713024857441768
$ws.Range("H14").Value = 3.461132737980023
$ws.Range("I14").Value = 0.6076242359217584
$ws.Range("J14").Value = 3.921643905914372
$ws.Range("F15").Value = 2.914031874512331
$ws.Range("H15").Value = 3.549364824109009
$ws.Range("I15").Value = 0.5071369742717423
$ws.Range("J15").Value = 3.859879938766589
$ws.Range("L15").Value = 1.775326769330522
$ws.Range("F16").Value = 2.71835923095414
$ws.Range("H16").Value = 3.103258029216448
$ws.Range("I16").Value = 0.6891476966100742
$ws.Range("J16").Value = 3.968643637828944
$ws.Range("L16").Value = 1.124482273704517
$ws.Range("F17").Value = 2.58070702680408
$ws.Range("H17").Value = 2.984728619664634
$ws.Range("I17").Value = 0.7129914068430308
$ws.Range("J17").Value = 4.055253480979948
$ws.Range("L17").Value = 0.6545325715981586
$ws.Range("F18").Value = 2.646447249450995
$ws.Range("H18").Value = 3.057552617977406
$ws.Range("I18").Value = 0.4074878344424646
$ws.Range("J18").Value = 4.031563077136536
$ws.Range("L18").Value = 0.9970497501096921
$ws.Range("F19").Value = 2.930790275888763
$ws.Range("H19").Value = 3.59160502022614
$ws.Range("I19").Value = 0.5619756823561046
$ws.Range("J19").Value = 3.92460231278308
$ws.Range("L19").Value = 2.041429437187643
$ws.Range("F20").Value = 3.051435999246952
$ws.Range("H20").Value = 3.353642042453706
$ws.Range("I20").Value = 0.4249044346260378
$ws.Range("J20").Value = 3.89541284571393
$ws.Range("L20").Value = 2.440782936632707
$ws.Range("F21").Value = 3.070495435058034
$ws.Range("H21").Value = 3.542000770356989
$ws.Range("I21").Value = 0.4268472660096853
$ws.Range("J21").Value = 3.901743553532211
$ws.Range("L21").Value = 2.538364855703463
$ws.Range("F22").Value = 2.984332605994403
$ws.Range("H22").Value = 3.569382153133348
$ws.Range("I22").Value = 0.5007390406270682
$ws.Range("J22").Value = 3.940570747341519
$ws.Range("L22").Value = 2.349968233186741
$ws.Range("F23").Value = 2.897713609847562
$ws.Range("H23").Value = 3.087841910245526
$ws.Range("I23").Value = 0.5291812806323573
$ws.Range("J23").Value = 3.977664966950252
$ws.Range("L23").Value = 2.157990508345986
$ws.Range("F24").Value = 3.471467837435296
$ws.Range("H24").Value = 3.960322399360249
$ws.Range("I24").Value = 0.6896507418452442
$ws.Range("J24").Value = 3.819206269810802
$ws.Range("L24").Value = 3.631962745272848
$ws.Range("F25").Value = 3.578279186158984
$ws.Range("H25").Value = 4.045898335052354
$ws.Range("I25").Value = 0.4235226668603078
$ws.Range("J25").Value = 3.805887039820029
$ws.Range("L25").Value = 3.880779623543889
